$wb = $excel.ActiveWorkbook

# "Statistics" sheet -- Production Date (E7) and Batch Number (E8)
$wsStats = $wb.Worksheets.Item("Statistics@0x5")
$wsStats.Range("E7").Formula = "'20191009"
$wsStats.Range("E8").Formula = "'1"

# "Calibration0" sheet -- Acceleration X K/D (E2/E3) and Voltage Battery D (E9)
$wsCal = $wb.Worksheets.Item("Calibration0@0x8")
$wsCal.Range("E2").Formula = "'0.0030518043786287308"
$wsCal.Range("E3").Formula = "'-99.62920379638672"
$wsCal.Range("E9").Formula = "'0.07954223453998566"
